$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("F").Insert()

$src = $ws.Range("C2:E20")
$dst = $ws.Range("D2:F20")
$dst.Value2 = $src.Value2
$ws.Range("C2:C20").Value2 = "Application"

$ws.Range("F1").Value2 = "field4"

$ws.Columns("D").ColumnWidth = 21.16
$ws.Columns("F").ColumnWidth = 10.16

Write-Host "Done"
